$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# ---- PIR sheet: append new rows ----
$ws = $wb.Worksheets.Item("PIR")
Set-TextCell $ws 284 1 '2026-02-04'
$ws.Cells.Item(284, 2).Value = '14:24:50'
$ws.Cells.Item(284, 3).Value = '14:00'
$ws.Cells.Item(284, 4).Value = 'Bathroom'
$ws.Cells.Item(284, 5).Value = 'No Motion'
$ws.Cells.Item(284, 6).Value = 'Inactive'
Set-TextCell $ws 285 1 '2026-02-04'
$ws.Cells.Item(285, 2).Value = '14:24:52'
$ws.Cells.Item(285, 3).Value = '14:00'
$ws.Cells.Item(285, 4).Value = 'Bathroom'
$ws.Cells.Item(285, 5).Value = 'No Motion'
$ws.Cells.Item(285, 6).Value = 'Inactive'
Set-TextCell $ws 286 1 '2026-02-04'
$ws.Cells.Item(286, 2).Value = '14:24:53'
$ws.Cells.Item(286, 3).Value = '14:00'
$ws.Cells.Item(286, 4).Value = 'Bathroom'
$ws.Cells.Item(286, 5).Value = 'Motion Detected'
$ws.Cells.Item(286, 6).Value = 'Active'
Set-TextCell $ws 287 1 '2026-02-04'
$ws.Cells.Item(287, 2).Value = '14:25:02'
$ws.Cells.Item(287, 3).Value = '14:00'
$ws.Cells.Item(287, 4).Value = 'Bathroom'
$ws.Cells.Item(287, 5).Value = 'No Motion'
$ws.Cells.Item(287, 6).Value = 'Inactive'
Set-TextCell $ws 288 1 '2026-02-04'
$ws.Cells.Item(288, 2).Value = '14:25:02'
$ws.Cells.Item(288, 3).Value = '14:00'
$ws.Cells.Item(288, 4).Value = 'Bathroom'
$ws.Cells.Item(288, 5).Value = 'Motion Detected'
$ws.Cells.Item(288, 6).Value = 'Active'
Set-TextCell $ws 289 1 '2026-02-04'
$ws.Cells.Item(289, 2).Value = '14:25:08'
$ws.Cells.Item(289, 3).Value = '14:00'
$ws.Cells.Item(289, 4).Value = 'Bathroom'
$ws.Cells.Item(289, 5).Value = 'No Motion'
$ws.Cells.Item(289, 6).Value = 'Inactive'
Set-TextCell $ws 290 1 '2026-02-04'
$ws.Cells.Item(290, 2).Value = '14:25:13'
$ws.Cells.Item(290, 3).Value = '14:00'
$ws.Cells.Item(290, 4).Value = 'Bathroom'
$ws.Cells.Item(290, 5).Value = 'No Motion'
$ws.Cells.Item(290, 6).Value = 'Inactive'
Set-TextCell $ws 291 1 '2026-02-04'
$ws.Cells.Item(291, 2).Value = '14:25:17'
$ws.Cells.Item(291, 3).Value = '14:00'
$ws.Cells.Item(291, 4).Value = 'Bathroom'
$ws.Cells.Item(291, 5).Value = 'Motion Detected'
$ws.Cells.Item(291, 6).Value = 'Active'
Set-TextCell $ws 292 1 '2026-02-04'
$ws.Cells.Item(292, 2).Value = '14:25:23'
$ws.Cells.Item(292, 3).Value = '14:00'
$ws.Cells.Item(292, 4).Value = 'Bathroom'
$ws.Cells.Item(292, 5).Value = 'No Motion'
$ws.Cells.Item(292, 6).Value = 'Inactive'
Set-TextCell $ws 293 1 '2026-02-04'
$ws.Cells.Item(293, 2).Value = '14:25:27'
$ws.Cells.Item(293, 3).Value = '14:00'
$ws.Cells.Item(293, 4).Value = 'Bathroom'
$ws.Cells.Item(293, 5).Value = 'Motion Detected'
$ws.Cells.Item(293, 6).Value = 'Active'
Set-TextCell $ws 294 1 '2026-02-04'
$ws.Cells.Item(294, 2).Value = '14:25:33'
$ws.Cells.Item(294, 3).Value = '14:00'
$ws.Cells.Item(294, 4).Value = 'Bathroom'
$ws.Cells.Item(294, 5).Value = 'No Motion'
$ws.Cells.Item(294, 6).Value = 'Inactive'
Set-TextCell $ws 295 1 '2026-02-04'
$ws.Cells.Item(295, 2).Value = '14:25:38'
$ws.Cells.Item(295, 3).Value = '14:00'
$ws.Cells.Item(295, 4).Value = 'Bathroom'
$ws.Cells.Item(295, 5).Value = 'No Motion'
$ws.Cells.Item(295, 6).Value = 'Inactive'
Set-TextCell $ws 296 1 '2026-02-04'
$ws.Cells.Item(296, 2).Value = '14:25:40'
$ws.Cells.Item(296, 3).Value = '14:00'
$ws.Cells.Item(296, 4).Value = 'Bathroom'
$ws.Cells.Item(296, 5).Value = 'Motion Detected'
$ws.Cells.Item(296, 6).Value = 'Active'

# ---- Humidity sheet: append new rows ----
$ws = $wb.Worksheets.Item("Humidity")
Set-TextCell $ws 237 1 '2026-02-04'
$ws.Cells.Item(237, 2).Value = '14:24:49'
$ws.Cells.Item(237, 3).Value = '14:00'
$ws.Cells.Item(237, 4).Value = 'Bathroom'
Set-TextCell $ws 237 5 '79.6%'
$ws.Cells.Item(237, 6).Value = 'Active'
Set-TextCell $ws 238 1 '2026-02-04'
$ws.Cells.Item(238, 2).Value = '14:24:51'
$ws.Cells.Item(238, 3).Value = '14:00'
$ws.Cells.Item(238, 4).Value = 'Bathroom'
Set-TextCell $ws 238 5 '78.8%'
$ws.Cells.Item(238, 6).Value = 'Active'
Set-TextCell $ws 239 1 '2026-02-04'
$ws.Cells.Item(239, 2).Value = '14:25:00'
$ws.Cells.Item(239, 3).Value = '14:00'
$ws.Cells.Item(239, 4).Value = 'Bathroom'
Set-TextCell $ws 239 5 '78.8%'
$ws.Cells.Item(239, 6).Value = 'Active'
Set-TextCell $ws 240 1 '2026-02-04'
$ws.Cells.Item(240, 2).Value = '14:25:05'
$ws.Cells.Item(240, 3).Value = '14:00'
$ws.Cells.Item(240, 4).Value = 'Bathroom'
Set-TextCell $ws 240 5 '79.7%'
$ws.Cells.Item(240, 6).Value = 'Active'
Set-TextCell $ws 241 1 '2026-02-04'
$ws.Cells.Item(241, 2).Value = '14:25:10'
$ws.Cells.Item(241, 3).Value = '14:00'
$ws.Cells.Item(241, 4).Value = 'Bathroom'
Set-TextCell $ws 241 5 '78.7%'
$ws.Cells.Item(241, 6).Value = 'Active'
Set-TextCell $ws 242 1 '2026-02-04'
$ws.Cells.Item(242, 2).Value = '14:25:15'
$ws.Cells.Item(242, 3).Value = '14:00'
$ws.Cells.Item(242, 4).Value = 'Bathroom'
Set-TextCell $ws 242 5 '79.7%'
$ws.Cells.Item(242, 6).Value = 'Active'
Set-TextCell $ws 243 1 '2026-02-04'
$ws.Cells.Item(243, 2).Value = '14:25:20'
$ws.Cells.Item(243, 3).Value = '14:00'
$ws.Cells.Item(243, 4).Value = 'Bathroom'
Set-TextCell $ws 243 5 '78.8%'
$ws.Cells.Item(243, 6).Value = 'Active'
Set-TextCell $ws 244 1 '2026-02-04'
$ws.Cells.Item(244, 2).Value = '14:25:25'
$ws.Cells.Item(244, 3).Value = '14:00'
$ws.Cells.Item(244, 4).Value = 'Bathroom'
Set-TextCell $ws 244 5 '79.7%'
$ws.Cells.Item(244, 6).Value = 'Active'
Set-TextCell $ws 245 1 '2026-02-04'
$ws.Cells.Item(245, 2).Value = '14:25:35'
$ws.Cells.Item(245, 3).Value = '14:00'
$ws.Cells.Item(245, 4).Value = 'Bathroom'
Set-TextCell $ws 245 5 '78.7%'
$ws.Cells.Item(245, 6).Value = 'Active'

# ---- Temperature sheet: append new rows ----
$ws = $wb.Worksheets.Item("Temperature")
Set-TextCell $ws 236 1 '2026-02-04'
$ws.Cells.Item(236, 2).Value = '14:24:48'
$ws.Cells.Item(236, 3).Value = '14:00'
$ws.Cells.Item(236, 4).Value = 'Bathroom'
$ws.Cells.Item(236, 5).Value = '24.3C'
$ws.Cells.Item(236, 6).Value = 'Active'
Set-TextCell $ws 237 1 '2026-02-04'
$ws.Cells.Item(237, 2).Value = '14:24:49'
$ws.Cells.Item(237, 3).Value = '14:00'
$ws.Cells.Item(237, 4).Value = 'Bathroom'
$ws.Cells.Item(237, 5).Value = '24.3C'
$ws.Cells.Item(237, 6).Value = 'Active'
Set-TextCell $ws 238 1 '2026-02-04'
$ws.Cells.Item(238, 2).Value = '14:24:52'
$ws.Cells.Item(238, 3).Value = '14:00'
$ws.Cells.Item(238, 4).Value = 'Bathroom'
$ws.Cells.Item(238, 5).Value = '24.3C'
$ws.Cells.Item(238, 6).Value = 'Active'
Set-TextCell $ws 239 1 '2026-02-04'
$ws.Cells.Item(239, 2).Value = '14:25:01'
$ws.Cells.Item(239, 3).Value = '14:00'
$ws.Cells.Item(239, 4).Value = 'Bathroom'
$ws.Cells.Item(239, 5).Value = '24.3C'
$ws.Cells.Item(239, 6).Value = 'Active'
Set-TextCell $ws 240 1 '2026-02-04'
$ws.Cells.Item(240, 2).Value = '14:25:06'
$ws.Cells.Item(240, 3).Value = '14:00'
$ws.Cells.Item(240, 4).Value = 'Bathroom'
$ws.Cells.Item(240, 5).Value = '24.3C'
$ws.Cells.Item(240, 6).Value = 'Active'
Set-TextCell $ws 241 1 '2026-02-04'
$ws.Cells.Item(241, 2).Value = '14:25:11'
$ws.Cells.Item(241, 3).Value = '14:00'
$ws.Cells.Item(241, 4).Value = 'Bathroom'
$ws.Cells.Item(241, 5).Value = '24.3C'
$ws.Cells.Item(241, 6).Value = 'Active'
Set-TextCell $ws 242 1 '2026-02-04'
$ws.Cells.Item(242, 2).Value = '14:25:16'
$ws.Cells.Item(242, 3).Value = '14:00'
$ws.Cells.Item(242, 4).Value = 'Bathroom'
$ws.Cells.Item(242, 5).Value = '24.3C'
$ws.Cells.Item(242, 6).Value = 'Active'
Set-TextCell $ws 243 1 '2026-02-04'
$ws.Cells.Item(243, 2).Value = '14:25:21'
$ws.Cells.Item(243, 3).Value = '14:00'
$ws.Cells.Item(243, 4).Value = 'Bathroom'
$ws.Cells.Item(243, 5).Value = '24.3C'
$ws.Cells.Item(243, 6).Value = 'Active'
Set-TextCell $ws 244 1 '2026-02-04'
$ws.Cells.Item(244, 2).Value = '14:25:26'
$ws.Cells.Item(244, 3).Value = '14:00'
$ws.Cells.Item(244, 4).Value = 'Bathroom'
$ws.Cells.Item(244, 5).Value = '24.3C'
$ws.Cells.Item(244, 6).Value = 'Active'
Set-TextCell $ws 245 1 '2026-02-04'
$ws.Cells.Item(245, 2).Value = '14:25:36'
$ws.Cells.Item(245, 3).Value = '14:00'
$ws.Cells.Item(245, 4).Value = 'Bathroom'
$ws.Cells.Item(245, 5).Value = '24.4C'
$ws.Cells.Item(245, 6).Value = 'Active'

